$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OrderCreationCases")

# Insert a new row above row 4 (shifts rows 4-7 down to 5-8)
$ws.Rows.Item(4).Insert()

# New row 4: Trigger Proforma Order with New Customer
$ws.Cells.Item(4, 1).Value = "QA2"
$ws.Cells.Item(4, 2).Value = "Trigger Proforma Order with New Customer"
$ws.Cells.Item(4, 3).Value = "GB"
$ws.Cells.Item(4, 4).Value = "Proforma"

# New row 9 (appended at the end): Trigger Proforma Order with Existing Customer
$ws.Cells.Item(9, 1).Value = "QA2"
$ws.Cells.Item(9, 2).Value = "Trigger Proforma Order with Existing Customer"
$ws.Cells.Item(9, 3).Value = "GB"
$ws.Cells.Item(9, 4).Value = "Proforma"
$ws.Cells.Item(9, 5).Value = "20240725122330Test"
$ws.Cells.Item(9, 6).Value = "20240725122330Auto"
$ws.Cells.Item(9, 7).Value = "20240725122330@Wiley.com"

$ws.Range("B11").Select()
